$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98; existing rows 98-145 shift down to 99-146.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record.
$ws.Cells.Item(98, 1).Value = 1
$ws.Cells.Item(98, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(98, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(98, 4).Value = 44489
$ws.Cells.Item(98, 5).Value = 15
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100108
$ws.Cells.Item(98, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(98, 9).Value = 100108006
$ws.Cells.Item(98, 10).Value = "Plátano"
$ws.Cells.Item(98, 11).Value = "Sin especificar"
$ws.Cells.Item(98, 12).Value = "Pintón"
$ws.Cells.Item(98, 13).Value = 120
$ws.Cells.Item(98, 14).Value = 24000
$ws.Cells.Item(98, 15).Value = 25000
$ws.Cells.Item(98, 16).Value = 24500
$ws.Cells.Item(98, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(98, 18).Value = "Ecuador"
$ws.Cells.Item(98, 19).Value = 1225
$ws.Cells.Item(98, 20).Value = 20
